# Update countries & provincias Spain
# Refresh COVID-19 country statistics table with the latest numbers.
# A handful of countries changed rank order (Bulgaria overtook Sudan,
# Paraguay overtook Albania), so the country name in those rows also
# changes along with the refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Agosto de 2020 a las 02:42"

# Helper to write a full data row: country name + B..H stats
function Set-Row($row, $country, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Row 4 - Estados Unidos (country unchanged, stats refreshed)
Set-Row 4 "Estados Unidos" 4705817 70832 2327530 2221543 0 1459 156744

# Row 25 - Canada (country unchanged, stats refreshed)
Set-Row 25 "Canada" 116312 513 101227 6150 0 6 8935

# Row 81 - now Bulgaria (new, higher figures; Bulgaria overtakes Sudan)
Set-Row 81 "Bulgaria" 11690 270 6319 4988 0 9 383

# Row 82 - now Sudan (keeps former row-81 figures, dropped one rank)
Set-Row 82 "Sudan" 11644 148 6119 4779 0 21 746

# Row 98 - now Paraguay (new, higher figures; Paraguay overtakes Albania)
Set-Row 98 "Paraguay" 5338 131 3548 1741 0 2 49

# Row 99 - now Albania (keeps former row-98 figures, dropped one rank)
Set-Row 99 "Albania" 5276 79 2952 2167 0 3 157

# Row 139 - Uruguay (country unchanged, stats refreshed)
Set-Row 139 "Uruguay" 1264 21 994 235 0 0 35

# Row 166 - Guyana (country unchanged, stats refreshed)
Set-Row 166 "Guyana" 413 12 185 208 0 0 20

# Row 175 - Guadalupe (country unchanged, stats refreshed)
Set-Row 175 "Guadalupe" 265 21 179 72 0 0 14
